# Add a dev-tunnels documentation link (as a new Heading1-styled paragraph)
# right after the "Web application working on local but not working on
# azure! Asp.Net MVC (c#)" paragraph, followed by a blank Heading1-styled
# paragraph, mirroring the style of the other headings in this doc.

$d = $word.ActiveDocument

# Locate the paragraph that ends in "... Asp.Net MVC (c#)" and, scoped to
# just that paragraph's range, turn its trailing ")" into ")" followed by
# two new paragraph breaks: one carrying the URL text, one left blank.
# Using Find/Replace (rather than Range.InsertParagraphAfter) keeps the
# inserted paragraphs free of stray empty runs, and both new paragraphs
# inherit the paragraph's own Heading1 formatting automatically.
$target = $d.Paragraphs(3).Range

$target.Find.Execute(
    ")",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    ")^phttps://learn.microsoft.com/en-us/aspnet/core/test/dev-tunnels?view=aspnetcore-8.0^p",
    2
) | Out-Null
